$wb = $excel.ActiveWorkbook

# --- Update the "Yearly" sheet ---
$yearly = $wb.Worksheets.Item("Yearly")
$yearly.Range("D13").Value = 102.41

# Update selection on Yearly sheet to J12
$yearly.Range("J12").Select()

# --- Update the "All Time" sheet ---
$allTime = $wb.Worksheets.Item("All Time")
$allTime.Activate()
$allTime.Range("I13").Select()
$excel.ActiveWindow.ScrollRow = 31
